$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Find.Execute("431×4=", $true, $false, $false, $false, $false, $true, 1, $false, "433×2=", 1) | Out-Null
$cell = $t.Cell(1, 2)
$cell.Range.Find.Execute("390×7=", $true, $false, $false, $false, $false, $true, 1, $false, "456×7=", 1) | Out-Null
$cell = $t.Cell(1, 3)
$cell.Range.Find.Execute("657×6=", $true, $false, $false, $false, $false, $true, 1, $false, "736×9=", 1) | Out-Null
$cell = $t.Cell(1, 4)
$cell.Range.Find.Execute("747×6=", $true, $false, $false, $false, $false, $true, 1, $false, "495×5=", 1) | Out-Null
$cell = $t.Cell(1, 5)
$cell.Range.Find.Execute("952×2=", $true, $false, $false, $false, $false, $true, 1, $false, "541×7=", 1) | Out-Null
$cell = $t.Cell(5, 1)
$cell.Range.Find.Execute("444×8=", $true, $false, $false, $false, $false, $true, 1, $false, "672×3=", 1) | Out-Null
$cell = $t.Cell(5, 2)
$cell.Range.Find.Execute("460×8=", $true, $false, $false, $false, $false, $true, 1, $false, "753×5=", 1) | Out-Null
$cell = $t.Cell(5, 3)
$cell.Range.Find.Execute("480×3=", $true, $false, $false, $false, $false, $true, 1, $false, "958×2=", 1) | Out-Null
$cell = $t.Cell(5, 4)
$cell.Range.Find.Execute("714×9=", $true, $false, $false, $false, $false, $true, 1, $false, "741×7=", 1) | Out-Null
$cell = $t.Cell(5, 5)
$cell.Range.Find.Execute("604×7=", $true, $false, $false, $false, $false, $true, 1, $false, "782×4=", 1) | Out-Null
$cell = $t.Cell(10, 1)
$cell.Range.Find.Execute("926×5=", $true, $false, $false, $false, $false, $true, 1, $false, "430×2=", 1) | Out-Null
$cell = $t.Cell(10, 2)
$cell.Range.Find.Execute("341×2=", $true, $false, $false, $false, $false, $true, 1, $false, "427×4=", 1) | Out-Null
$cell = $t.Cell(10, 3)
$cell.Range.Find.Execute("948×6=", $true, $false, $false, $false, $false, $true, 1, $false, "433×2=", 1) | Out-Null
$cell = $t.Cell(10, 4)
$cell.Range.Find.Execute("421×5=", $true, $false, $false, $false, $false, $true, 1, $false, "257×8=", 1) | Out-Null
$cell = $t.Cell(10, 5)
$cell.Range.Find.Execute("451×7=", $true, $false, $false, $false, $false, $true, 1, $false, "497×4=", 1) | Out-Null
$cell = $t.Cell(15, 1)
$cell.Range.Find.Execute("970×5=", $true, $false, $false, $false, $false, $true, 1, $false, "891×9=", 1) | Out-Null
$cell = $t.Cell(15, 2)
$cell.Range.Find.Execute("555×3=", $true, $false, $false, $false, $false, $true, 1, $false, "616×8=", 1) | Out-Null
$cell = $t.Cell(15, 3)
$cell.Range.Find.Execute("429×9=", $true, $false, $false, $false, $false, $true, 1, $false, "238×9=", 1) | Out-Null
$cell = $t.Cell(15, 4)
$cell.Range.Find.Execute("390×7=", $true, $false, $false, $false, $false, $true, 1, $false, "571×6=", 1) | Out-Null
$cell = $t.Cell(15, 5)
$cell.Range.Find.Execute("593×8=", $true, $false, $false, $false, $false, $true, 1, $false, "182×9=", 1) | Out-Null
$cell = $t.Cell(20, 1)
$cell.Range.Find.Execute("866×9=", $true, $false, $false, $false, $false, $true, 1, $false, "885×9=", 1) | Out-Null
$cell = $t.Cell(20, 2)
$cell.Range.Find.Execute("688×6=", $true, $false, $false, $false, $false, $true, 1, $false, "505×6=", 1) | Out-Null
$cell = $t.Cell(20, 3)
$cell.Range.Find.Execute("750×5=", $true, $false, $false, $false, $false, $true, 1, $false, "707×5=", 1) | Out-Null
$cell = $t.Cell(20, 4)
$cell.Range.Find.Execute("645×6=", $true, $false, $false, $false, $false, $true, 1, $false, "540×9=", 1) | Out-Null
$cell = $t.Cell(20, 5)
$cell.Range.Find.Execute("661×7=", $true, $false, $false, $false, $false, $true, 1, $false, "718×9=", 1) | Out-Null
